# Weekly fruit/vegetable price update: insert one new observation row
# right before the existing row 216, shifting all subsequent rows down
# by one (rows 216-345 become 217-346), and fill the newly inserted
# row 216 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 216 (pushes old 216..345 -> 217..346,
# carrying along their existing formatting/number formats).
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new record.
$ws.Cells.Item(216, 1).Value  = 4
$ws.Cells.Item(216, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(216, 3).Value  = "Los Lagos"
$ws.Cells.Item(216, 4).Value  = 44879
$ws.Cells.Item(216, 5).Value  = 10
$ws.Cells.Item(216, 6).Value  = 100112017
$ws.Cells.Item(216, 7).Value  = "Apio"
$ws.Cells.Item(216, 8).Value  = "Americana (o)"
$ws.Cells.Item(216, 9).Value  = "Primera"
$ws.Cells.Item(216, 10).Value = 20
$ws.Cells.Item(216, 11).Value = 15000
$ws.Cells.Item(216, 12).Value = 15000
$ws.Cells.Item(216, 13).Value = 15000
$ws.Cells.Item(216, 14).Value = "$/docena de matas"
$ws.Cells.Item(216, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(216, 16).Value = 2500
$ws.Cells.Item(216, 17).Value = 6
$ws.Cells.Item(216, 18).Value = "Hortaliza"
